$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so values like "63.746.30" or "+3.33%" are not
# reinterpreted as numbers/dates by Excel.
$cells = @( "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "D14", "E14", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "E28", "E29", "E30", "D31", "E31", "E32", "D33", "E33", "D34", "E34", "E35", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "E49", "B50", "C50", "D50", "E50", "B51", "C51", "D51", "E51" )
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.746.30"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "3.129.16"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "590.52"
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("D6").Value = "146.52"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.121.64"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +18.61%  "
$ws.Range("D11").Value = "5.68"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("E13").Value = "  +6.51%  "
$ws.Range("D14").Value = "35.97"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "3.649.07"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").Value = "63.696.99"
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "3.129.87"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("D20").Value = "465.11"
$ws.Range("E20").Value = "  +3.87%  "
$ws.Range("D21").Value = "14.19"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").Value = "0.732"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("D24").Value = "13.27"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").Value = "82.28"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "8.60"
$ws.Range("E27").Value = "  +7.34%  "
$ws.Range("E28").Value = "  +2.93%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("D31").Value = "6.83"
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "0.0₃0867"
$ws.Range("E34").Value = "  +7.14%  "
$ws.Range("E35").Value = "  +10.33%  "
$ws.Range("E36").Value = "  +2.59%  "
$ws.Range("D37").Value = "3.36"
$ws.Range("E37").Value = "  +13.19%  "
$ws.Range("D38").Value = "6.12"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").Value = "50.76"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").Value = "449.18"
$ws.Range("E40").Value = "  +7.87%  "
$ws.Range("D41").Value = "8.70"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "0.0372"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").Value = "2.919.41"
$ws.Range("E43").Value = "  +5.32%  "
$ws.Range("D44").Value = "0.278"
$ws.Range("E44").Value = "  +5.70%  "
$ws.Range("D45").Value = "0.111"
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("D47").Value = "127.52"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "34.15"
$ws.Range("E50").Value = "  -7.57%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "24.64"
$ws.Range("E51").Value = "  +2.39%  "
